$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8 (pushing existing rows 8..23 down to 9..24).
# FOURVENT SYRUP 125ML is alphabetically between "DURICEF..." (row 7) and
# "INJECTMOL..." (old row 8, now row 9).
$ws.Rows(8).Insert()

# Fill in the new row's values, following the same layout as the other data rows:
#   A = sequence number, B:G = product name (merged), H:K = balance ratio (merged),
#   L:M = price (merged), N = transactions ratio
$ws.Cells.Item(8, 1).Value = 5
$ws.Cells.Item(8, 2).Value = "FOURVENT SYRUP 125ML"
$ws.Cells.Item(8, 8).Value = "10:0"
$ws.Cells.Item(8, 12).Value = 24
$ws.Cells.Item(8, 14).Value = "1:0"

# Make sure the merged regions match the pattern used by every other data row.
$ws.Range("B8:G8").Merge()
$ws.Range("H8:K8").Merge()
$ws.Range("L8:M8").Merge()

$wb.Save()
